# Applies the "nuevos experimentos no convexos" edit:
# updates numeric/expression values across several sheets while
# keeping labels/headers intact. Most of the updated cells hold
# numeric-looking text (stored as shared strings in the original
# workbook), so we briefly force a text number format before writing
# the value and then restore the default style - this keeps Excel
# from silently re-typing them as real numbers.

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# NOTE: this workbook has two sheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively, so address those two sheets by position instead.

# --- Sheet: Restricciones_del_follower (position 3) --------------------
$ws = $wb.Worksheets.Item(3)

# Row 2 (J_0_L0_v)
Set-TextValue $ws.Range("A2") "1.4164683088845837 - 0.4882369222252976y_1 + 0.49986161084970937y_2"
Set-TextValue $ws.Range("B2") "-1.4164683088845837"
Set-TextValue $ws.Range("C2") "J_0_L0_v"
Set-TextValue $ws.Range("D2") "0.35"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "0"

# Row 3 (J_0_L0_v)
Set-TextValue $ws.Range("A3") "-1.4164683088845837 + 0.4882369222252976y_1 - 0.49986161084970937y_2"
Set-TextValue $ws.Range("B3") "-2.5835316911154163"
Set-TextValue $ws.Range("C3") "J_0_L0_v"
Set-TextValue $ws.Range("D3") "0.91"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "0"

# Row 4 (J_0_LP_v)
Set-TextValue $ws.Range("A4") "20.584306670357044 - 2x - 1.51120952117354y_1 + 1.5471907002491005y_2"
Set-TextValue $ws.Range("B4") "-36.584306670357044"
Set-TextValue $ws.Range("C4") "J_0_LP_v"
Set-TextValue $ws.Range("D4") "0.78"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "0"

# Row 5 (J_Ne_L0_v)
Set-TextValue $ws.Range("A5") "-67.01646830888458 + 8x + 0.4882369222252976y_1 - 0.49986161084970937y_2"
Set-TextValue $ws.Range("B5") "18.21646830888458"
Set-TextValue $ws.Range("C5") "J_Ne_L0_v"
Set-TextValue $ws.Range("D5") "0.57"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "0"

# Row 6 (J_Ne_L0_v)
Set-TextValue $ws.Range("A6") "4.967063382230833 - 2x - 0.9764738444505952y_1 + 0.9997232216994187y_2"
Set-TextValue $ws.Range("B6") "-7.032936617769167"
Set-TextValue $ws.Range("C6") "J_Ne_L0_v"
Set-TextValue $ws.Range("D6") "0.0"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "0"

# --- Sheet: Punto_modificado (position 4) -------------------------------
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "8.1"
Set-TextValue $ws.Range("B2") "5.0"
Set-TextValue $ws.Range("C2") "2.05"

# --- Sheet: Vector_bf (position 5) --------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "1.6270357044007748"
Set-TextValue $ws.Range("A3") "-0.6419651259341269"

# --- Sheet: Vector_BF (position 6) --------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "1.0"
Set-TextValue $ws.Range("A3") "3.0"
Set-TextValue $ws.Range("A4") "-2.0"

# --- Sheet: Vector_Alpha (position 7) -----------------------------------
# A2/A3 are stored as real numbers (not shared strings) in this sheet.
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.29
$ws.Range("A3").Value = 1.26
